$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking Price cells as Text so Excel doesn't
# auto-convert them to numbers (which would strip the exact digits,
# e.g. trailing zeros) - matches the inline-string cell type in the source.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '25.844.63'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.741.36'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '225.35'
$ws.Range('E5').Value = '  -4.96%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '0.5161'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('D8').Value = '0.2792'
$ws.Range('E8').Value = '  +6.40%  '
$ws.Range('D9').Value = '39.10'
$ws.Range('E9').Value = '  -5.10%  '
$ws.Range('D10').Value = '0.06094'
$ws.Range('D11').Value = '1.742.00'
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').Value = '0.06973'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '15.20'
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').Value = '0.6349'
$ws.Range('E14').Value = '  +5.34%  '
$ws.Range('D15').Value = '4.497'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').Value = '76.47'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '25.864.64'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '11.43'
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').Value = '0.000006582'
$ws.Range('E21').Value = '  -3.37%  '
$ws.Range('D22').Value = '1.962.53'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D24').Value = '8.429'
$ws.Range('D25').Value = '5.097'
$ws.Range('E25').Value = '  -1.32%  '
$ws.Range('D26').Value = '137.56'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  +3.33%  '
$ws.Range('D28').Value = '1.817'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D29').Value = '14.96'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').Value = '102.51'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').Value = '0.08273'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('D32').Value = '3.610'
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('D33').Value = '3.408'
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').Value = '0.04399'
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('D35').Value = '2.620'
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('D36').Value = '0.9707'
$ws.Range('E36').Value = '  -3.09%  '
$ws.Range('D37').Value = '0.5998'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -2.40%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = '0.9996'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').Value = '1.899'
$ws.Range('E41').Value = '  -1.34%  '
$ws.Range('D42').Value = '100.73'
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').Value = '0.3821'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '0.7242'
$ws.Range('E44').Value = '  -2.75%  '
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').Value = '0.05458'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('D47').Value = '6.262'
$ws.Range('E47').Value = '  +5.46%  '
$ws.Range('D48').Value = '0.1099'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('D49').Value = '29.73'
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('D50').Value = '52.09'
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = '7.454'
$ws.Range('E51').Value = '  -0.95%  '

# Restore the default cell style (the source cells carry no explicit
# style index) now that the values are safely stored as text.
$ws.Range('D2:D51').Style = 'Normal'
